$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-theme the deck: swap the slide master's theme palette from the old
#    "Integral" design colours over to the standard Office theme colours
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), matching the new design that
#    was applied to the presentation.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# index : name     : new RGB (hex)
$newColors = @(
    @(1,  0x000000),  # dk1
    @(2,  0xFFFFFF),  # lt1
    @(3,  0x44546A),  # dk2
    @(4,  0xE7E6E6),  # lt2
    @(5,  0x5B9BD5),  # accent1
    @(6,  0xED7D31),  # accent2
    @(7,  0xA5A5A5),  # accent3
    @(8,  0xFFC000),  # accent4
    @(9,  0x4472C4),  # accent5
    @(10, 0x70AD47),  # accent6
    @(11, 0x0563C1),  # hlink
    @(12, 0x954F72)   # folHlink
)

foreach ($entry in $newColors) {
    $idx = $entry[0]
    $hex = $entry[1]
    # $hex is written as 0xRRGGBB; pull out the three channel bytes ...
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # ... then repack as the little-endian 0x00BBGGRR long that the
    # ColorFormat.RGB property (and VBA's RGB() helper) expects.
    $bgr = $r + ($g * 256) + ($b * 65536)
    $themeColors.Colors($idx).RGB = $bgr
}

# ---------------------------------------------------------------------------
# 2) Point the summary table on slide 16 at the new built-in table style.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{A958BEF4-563C-40CD-9898-CCFB3343E13A}")
    }
}
